# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Fri Jul 14 17:25:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (values like "31.334.23" or "1.007"
# use "." as a separator and are NOT real numbers). Force those ranges to
# Text format before writing so Excel keeps the new prices as literal text
# instead of silently re-parsing them as numbers (which would drop trailing
# zeros / merge thousand separators).
$ws.Range("D2:D37").NumberFormat = "@"
$ws.Range("D39:D46").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

$ws.Range("D2").Value = '31.158.95'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.986.53'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +1.12%  '
$ws.Range("D5").Value = '254.35'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '0.7285'
$ws.Range("E6").Value = '  +15.46%  '
$ws.Range("D7").Value = '1.008'
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("D8").Value = '0.3380'
$ws.Range("E8").Value = '  +4.26%  '
$ws.Range("D9").Value = '27.38'
$ws.Range("E9").Value = '  +8.52%  '
$ws.Range("D10").Value = '0.07112'
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("D11").Value = '0.8274'
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").Value = '0.08115'
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("D13").Value = '1.992.89'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.572'
$ws.Range("E14").Value = '  +3.81%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '98.95'
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("D16").Value = '15.28'
$ws.Range("E16").Value = '  +10.94%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").Value = '267.24'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '31.203.54'
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").Value = '6.091'
$ws.Range("E19").Value = '  +8.28%  '
$ws.Range("D20").Value = '0.000008210'
$ws.Range("E20").Value = '  +6.73%  '
$ws.Range("D21").Value = '2.260.71'
$ws.Range("E21").Value = '  +3.26%  '
$ws.Range("D22").Value = '1.009'
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '1.011'
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("D24").Value = '7.041'
$ws.Range("E24").Value = '  +6.08%  '
$ws.Range("D25").Value = '9.993'
$ws.Range("E25").Value = '  +5.48%  '
$ws.Range("D26").Value = '162.36'
$ws.Range("E26").Value = '  -1.66%  '
$ws.Range("D27").Value = '19.69'
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("D28").Value = '2.350'
$ws.Range("E28").Value = '  +10.66%  '
$ws.Range("D29").Value = '0.1316'
$ws.Range("E29").Value = '  +7.80%  '
$ws.Range("D30").Value = '1.602'
$ws.Range("E30").Value = '  +3.10%  '
$ws.Range("D31").Value = '1.383'
$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("D32").Value = '4.608'
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").Value = '4.410'
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("D34").Value = '0.05326'
$ws.Range("E34").Value = '  +5.92%  '
$ws.Range("D35").Value = '1.285'
$ws.Range("E35").Value = '  +7.15%  '
$ws.Range("D36").Value = '0.7826'
$ws.Range("E36").Value = '  +7.33%  '
$ws.Range("D37").Value = '2.808'
$ws.Range("E37").Value = '  +3.50%  '
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("D39").Value = '2.894'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '83.55'
$ws.Range("E40").Value = '  +7.79%  '
$ws.Range("D41").Value = '6.777'
$ws.Range("E41").Value = '  +4.12%  '
$ws.Range("D42").Value = '0.4621'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("D43").Value = '2.102'
$ws.Range("E43").Value = '  +3.54%  '
$ws.Range("D44").Value = '0.8553'
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.008'
$ws.Range("E45").Value = '  +0.96%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '104.68'
$ws.Range("E46").Value = '  +2.02%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").Value = '7.621'
$ws.Range("E48").Value = '  +4.26%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.575'
$ws.Range("E49").Value = '  +10.97%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '37.21'
$ws.Range("E50").Value = '  +3.50%  '
$ws.Range("D51").Value = '0.4300'
$ws.Range("E51").Value = '  +2.68%  '
